# Swap the ASV_ID / Species_name / Common_name / Category values between
# rows 42 and 43, and move the numeric value that was in J43 (0) up to J42,
# clearing J43 since it becomes empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42 becomes what used to be row 43's identity columns.
$ws.Range("A42").Value = "c0a3f3ed23f04247d92740a9502f8b57"
$ws.Range("B42").Value = "unassigned"
$ws.Range("C42").Value = "unassigned"
$ws.Range("D42").Value = "unassigned"

# Row 43 becomes what used to be row 42's identity columns.
$ws.Range("A43").Value = "307c55294ffe3b8aa46fce358d55590e"
$ws.Range("B43").Value = "Homo sapiens"
$ws.Range("C43").Value = "Human"
$ws.Range("D43").Value = "Human"

# J43 (0) moves to J42; J43 becomes empty.
$ws.Range("J42").Value = 0
$ws.Range("J43").ClearContents()
